$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.986.89'
$ws.Range("E2").Value = '  +3.76%  '

$ws.Range("D3").Value = '1.678.04'
$ws.Range("E3").Value = '  +2.90%  '

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = "'219.83"
$ws.Range("E5").Value = '  +2.37%  '

$ws.Range("D6").Value = "'0.534"
$ws.Range("E6").Value = '  +2.38%  '

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").Value = "'29.19"
$ws.Range("E8").Value = '  +1.66%  '

$ws.Range("E9").Value = '  +2.73%  '

$ws.Range("D10").Value = "'0.0640"
$ws.Range("E10").Value = '  +5.17%  '

$ws.Range("D11").Value = "'0.0909"
$ws.Range("E11").Value = '  +0.94%  '

$ws.Range("D12").Value = '1.918.99'
$ws.Range("E12").Value = '  +2.94%  '

$ws.Range("D13").Value = '1.672.37'
$ws.Range("E13").Value = '  +2.65%  '

$ws.Range("D14").Value = "'10.14"
$ws.Range("E14").Value = '  +7.92%  '

$ws.Range("D15").Value = "'0.605"
$ws.Range("E15").Value = '  +6.69%  '

$ws.Range("D16").Value = "'4.10"
$ws.Range("E16").Value = '  +6.89%  '

$ws.Range("D17").Value = '30.907.33'
$ws.Range("E17").Value = '  +3.39%  '

$ws.Range("D18").Value = "'66.21"
$ws.Range("E18").Value = '  +1.15%  '

$ws.Range("D19").Value = "'247.22"
$ws.Range("E19").Value = '  +2.70%  '

$ws.Range("E20").Value = '  +2.39%  '

$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = '  -0.07%  '

$ws.Range("D22").Value = "'4.26"
$ws.Range("E22").Value = '  +3.13%  '

$ws.Range("D23").Value = "'10.00"
$ws.Range("E23").Value = '  +1.95%  '

$ws.Range("E24").Value = '  -0.56%  '

$ws.Range("D25").Value = "'158.92"
$ws.Range("E25").Value = '  +0.85%  '

$ws.Range("D26").Value = "'15.87"
$ws.Range("E26").Value = '  +2.48%  '

$ws.Range("E27").Value = '  +2.30%  '

$ws.Range("D28").Value = "'6.69"
$ws.Range("E28").Value = '  +1.60%  '

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("D30").Value = "'0.0495"
$ws.Range("E30").Value = '  +1.07%  '

$ws.Range("D31").Value = "'3.51"
$ws.Range("E31").Value = '  +3.89%  '

$ws.Range("E32").Value = '  +3.65%  '

$ws.Range("D33").Value = "'3.34"
$ws.Range("E33").Value = '  +5.04%  '

$ws.Range("D34").Value = '1.519.37'
$ws.Range("E34").Value = '  +6.52%  '

$ws.Range("E35").Value = '  +3.46%  '

$ws.Range("D36").Value = "'84.41"
$ws.Range("E36").Value = '  +12.60%  '

$ws.Range("E37").Value = '  +0.49%  '

$ws.Range("D38").Value = "'0.609"
$ws.Range("E38").Value = '  +9.45%  '

$ws.Range("D39").Value = "'0.0180"
$ws.Range("E39").Value = '  +5.26%  '

$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = "'2.68"
$ws.Range("E40").Value = '  -3.26%  '

$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D41").Value = "'2.29"
$ws.Range("E41").Value = '  +0.23%  '

$ws.Range("D42").Value = "'2.06"
$ws.Range("E42").Value = '  +3.78%  '

$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = "'0.839"
$ws.Range("E43").Value = '  +0.76%  '

$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").Value = "'0.0503"
$ws.Range("E44").Value = '  +0.66%  '

$ws.Range("E45").Value = '  +2.26%  '

$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").Value = "'5.58"
$ws.Range("E47").Value = '  +4.69%  '

$ws.Range("D48").Value = "'51.23"
$ws.Range("E48").Value = '  +4.65%  '

$ws.Range("D49").Value = '1.811.39'
$ws.Range("E49").Value = '  +2.25%  '

$ws.Range("D50").Value = '0.0₆0119'
$ws.Range("E50").Value = '  +8.25%  '

$ws.Range("D51").Value = "'93.34"
$ws.Range("E51").Value = '  +1.17%  '
